# Insert a new weekly price record as the new row 33, pushing all
# subsequent rows down by one (old row 33 becomes row 34, ..., old row 82
# becomes row 83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33:82 down to 34:83 by inserting a blank row at 33.
$ws.Rows("33:33").Insert()

# Populate the new record in row 33.
$ws.Range("A33").Value() = 1
$ws.Range("B33").Value() = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value() = "Arica y Parinacota"
$ws.Range("D33").Value() = 44810
$ws.Range("E33").Value() = 15
$ws.Range("F33").Value() = 100114001
$ws.Range("G33").Value() = "Papa"
$ws.Range("H33").Value() = "Asterix"
$ws.Range("I33").Value() = "1a (cosecha)"
$ws.Range("J33").Value() = 1000
$ws.Range("K33").Value() = 11000
$ws.Range("L33").Value() = 12000
$ws.Range("M33").Value() = 11500
$ws.Range("N33").Value() = "`$/saco 25 kilos"
$ws.Range("O33").Value() = "Región del Maule"
$ws.Range("P33").Value() = 460
$ws.Range("Q33").Value() = 25
$ws.Range("R33").Value() = "Hortaliza"
